# Sprint backlog da semana e tutorial nodeJS
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new sheet "Sprint Backlog (19-04)" as a copy of the most
#    recent sprint sheet ("Sprint Backlog (09-04)") so it inherits the
#    same layout/column widths/base formatting, then place it right
#    after that sheet (i.e. at the end of the tab strip).
# ---------------------------------------------------------------------
$wsPrev = $wb.Worksheets.Item("Sprint Backlog (09-04)")
$wsPrev.Copy($null, $wsPrev)
$ws7 = $wb.Worksheets.Item($wsPrev.Index + 1)
$ws7.Name = "Sprint Backlog (19-04)"

# ---------------------------------------------------------------------
# 2. Remove the stray note in sheet "Sprint Backlog (09-04)" (E7 = "OK")
#    that does not belong to the table any more.
# ---------------------------------------------------------------------
$wsPrev.Range("E7").ClearContents()
$wsPrev.Range("A1:D8").Select()

# ---------------------------------------------------------------------
# 3. Rebuild the merged "Grupo" column on the new sheet: instead of one
#    group spanning all 6 task rows, split it into three groups.
# ---------------------------------------------------------------------
$ws7.Range("B2:B7").UnMerge()
$ws7.Range("B2:B4").Merge()
$ws7.Range("B5:B6").Merge()
$ws7.Range("B7:B8").Merge()

# ---------------------------------------------------------------------
# 4. Populate the new sprint's tasks.
# ---------------------------------------------------------------------
$ws7.Range("B2").Value2 = "Desenvolvimento"
$ws7.Range("C2").Value2 = "Fazer barra de menu lateral para navegação entre as páginas do site."
$ws7.Range("D2").Value2 = "Essencial"

$ws7.Range("C3").Value2 = "Fazer modal com opções de edição de produto, ambiente e informações dos galpões na página galpões."
$ws7.Range("D3").Value2 = "Essencial"

$ws7.Range("C4").Value2 = "Alinhar cores das páginas já feitas e as conectar por links."
$ws7.Range("D4").Value2 = "Desejavel"

$ws7.Range("B5").Value2 = "Banco de dados"
$ws7.Range("C5").Value2 = "Inserir dados consistentes no banco de dados relacionados a usuários e galpões."
$ws7.Range("D5").Value2 = "Desejavel"

$ws7.Range("C6").Value2 = "Testar querys de select com as informações inseridas no banco de dados."
$ws7.Range("D6").Value2 = "Desejável"

$ws7.Range("B7").Value2 = "Documentação"
$ws7.Range("C7").Value2 = "Fazer tutorial de nodeJS para professora Marise."
$ws7.Range("D7").Value2 = "Desejavel"

$ws7.Range("C8").Value2 = "Revisar planilha de riscos do projeto."
$ws7.Range("D8").Value2 = "Desejável"

# ---------------------------------------------------------------------
# 5. Row heights on the new sheet grew because of the longer task
#    descriptions (wrapped text).
# ---------------------------------------------------------------------
$ws7.Rows.Item(2).RowHeight = 60
$ws7.Rows.Item(3).RowHeight = 105
$ws7.Rows.Item(4).RowHeight = 60.75
$ws7.Rows.Item(5).RowHeight = 75
$ws7.Rows.Item(6).RowHeight = 75.75
$ws7.Rows.Item(7).RowHeight = 45
$ws7.Rows.Item(8).RowHeight = 30.75

# ---------------------------------------------------------------------
# 6. View state: new sheet becomes the selected/active tab, zoomed in,
#    with a fresh selection.
# ---------------------------------------------------------------------
$ws7.Range("F5").Select()
$excel.ActiveWindow.Zoom = 70
$ws7.Activate()

$wb.Save()
